$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "50÷7="
$table.Cell(1, 2).Range.Text = "77÷4="
$table.Cell(1, 3).Range.Text = "98÷9="
$table.Cell(1, 4).Range.Text = "83÷4="
$table.Cell(1, 5).Range.Text = "41÷8="
$table.Cell(5, 1).Range.Text = "80÷2="
$table.Cell(5, 2).Range.Text = "21÷9="
$table.Cell(5, 3).Range.Text = "38÷3="
$table.Cell(5, 4).Range.Text = "41÷7="
$table.Cell(5, 5).Range.Text = "88÷8="
$table.Cell(9, 1).Range.Text = "23÷5="
$table.Cell(9, 2).Range.Text = "93÷4="
$table.Cell(9, 3).Range.Text = "51÷8="
$table.Cell(9, 4).Range.Text = "14÷9="
$table.Cell(9, 5).Range.Text = "59÷8="
$table.Cell(13, 1).Range.Text = "74÷7="
$table.Cell(13, 2).Range.Text = "30÷6="
$table.Cell(13, 3).Range.Text = "31÷5="
$table.Cell(13, 4).Range.Text = "85÷5="
$table.Cell(13, 5).Range.Text = "16÷2="
$table.Cell(17, 1).Range.Text = "54÷6="
$table.Cell(17, 2).Range.Text = "93÷8="
$table.Cell(17, 3).Range.Text = "47÷8="
$table.Cell(17, 4).Range.Text = "50÷3="
$table.Cell(17, 5).Range.Text = "57÷9="
